$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- register the underlined-Calibri font used by the rich-text note below ---
# (Excel registers a "known font" in styles.xml the first time underline is
# applied anywhere; doing it on a scratch cell far away and then clearing the
# cell keeps the sheet's used range untouched while still getting the font
# entry written out.)
$scratch = $ws.Range("Z100")
$scratch.Value = "x"
$scratch.Font.Underline = $true
$scratch.Clear() | Out-Null

# --- header note additions ---
$ws.Range("A4").Value = "lidar_values_and_settings.py settings:"
$ws.Range("S8").Value = "The path fields accept relative and absolute paths, but only in unix syntax"

# --- rename "scale x" -> "scale_x" ---
$ws.Range("A33").Value = "scale_x"

# --- widen column B a bit ---
$ws.Columns.Item(2).ColumnWidth = 17

# --- new hidden spacer rows ---
$ws.Rows.Item(26).Hidden = $true
$ws.Rows.Item(27).Hidden = $true
$ws.Rows.Item(28).Hidden = $true
$ws.Rows.Item(41).Hidden = $true
$ws.Rows.Item(42).Hidden = $true
$ws.Rows.Item(43).Hidden = $true

# --- new "read_pcap_from_file.py settings" block ---
$ws.Range("A45").Value = "read_pcap_from_file.py settings:"

$ws.Range("A46").Value = "From pcap"
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D46").Value = "If the path fields are left empty, or the path dir or file doesn't exist, the program will prompt you for manual selection"

$ws.Range("A47").Value = "flight_scans"
$ws.Range("B47").Value = "./flight_scans/wireshark_flight_1.pcap"

$ws.Range("A48").Value = "flight_ins"
$ws.Range("B48").Value = "./flight_ins/export_flight01_20.08.19_all.txt"

$ws.Range("A49").Value = "temp_location_frame_files"
$ws.Range("B49").Value = "./files_from_pcap/"
$ws.Range("D49").Value = "The location for saving the las files containing only one frame each"

$ws.Range("A51").Value = "packet_devisor"
$ws.Range("B51").Value = 1
$ws.Range("D51").Value = "For quick analysis you can choose to skip every n packet from being processed"

$ws.Range("A52").Value = "num_frames"
$ws.Range("B52").Value = 200
$ws.Range("D52").Value = "number of frames to be created from the pcap (it will read in the first n files)"
$ws.Range("D52").Characters(67, 5).Font.Underline = $true
$ws.Range("D52").Characters(72, 8).Font.Underline = $false
$ws.Range("D52").Characters(80, 1).Font.Underline = $true

$ws.Range("A54").Value = "From frame files"
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A55").Value = "num_frames_per_las_file"
$ws.Range("B55").Formula = "=B51*200"
$ws.Range("D55").Value = "The amount of frames to load and put together in one las file, keep in mind that this amount must be stored in memory"

$ws.Range("A56").Value = "outfile_directory"
$ws.Range("B56").Value = "./processed_las/garbage/"

# --- move the saved selection cursor, like the original author's session ---
$ws.Range("Q35").Select() | Out-Null
